# Applies the "new agri files and 18_CT bad values in production" edit:
#   - Inserts a new sub3sector value "09_01_02_gas_power_h2" into column G,
#     immediately after "09_01_02_gas_power_ccs" (before "09_01_03_oil"),
#     shifting the remaining G-column list down by one row.
#   - Inserts a new sub3sector value "18_01_02_gas_power_h2" into column G,
#     immediately after "18_01_02_gas_power_ccs" (before "18_01_03_oil"),
#     shifting the remaining G-column list down by one row.
#   - Moves the review comments that were anchored to the old row numbers so
#     they stay attached to the same logical G-column entries:
#       G20 ("should add 05_other")                    -> G21
#       G24 ("we should add 9_x_others, ...")           -> G25
#       G46 ("alex might think about own use in ccs...")-> G47

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 7
$lastRow  = 105

# --- Read the existing column G list (rows 7..105) ------------------------
$values = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $values += $ws.Cells.Item($r, 7).Value2
}

# --- Build the new list with the two insertions ----------------------------
$newValues = @()
foreach ($v in $values) {
    if ($v -eq "09_01_03_oil") {
        $newValues += "09_01_02_gas_power_h2"
    }
    if ($v -eq "18_01_03_oil") {
        $newValues += "18_01_02_gas_power_h2"
    }
    $newValues += $v
}

# --- Write the new list back out, now spanning rows 7..107 -----------------
$r = $firstRow
foreach ($v in $newValues) {
    $ws.Cells.Item($r, 7).Value2 = $v
    $r = $r + 1
}

# --- Re-anchor the comments that sit on the shifted G-column entries -------
function Move-CellComment($ws, $fromAddr, $toAddr) {
    $src = $ws.Range($fromAddr)
    $c = $src.Comment
    if ($c -eq $null) { return }
    $txt = $c.Text()
    $author = $c.Author
    $visible = $c.Visible
    $c.Delete()
    $dst = $ws.Range($toAddr)
    $newc = $dst.AddComment($txt)
    $newc.Author = $author
    $newc.Visible = $visible
}

Move-CellComment $ws "G46" "G47"
Move-CellComment $ws "G24" "G25"
Move-CellComment $ws "G20" "G21"

# --- Keep the view/dimension metadata in line with the new data extent -----
$ws.Range("G81").Select() | Out-Null
